# Add "DONE" / "PARTLY DONE" status markers to the front of specific
# checklist items in the "Worms detection algorithm" / "Evaluation" list.
#
# Colors (Word BGR-ordered decimal values for the OOXML RRGGBB hex):
#   00CC00 (green)  -> 52224
#   FF9900 (orange) -> 39423

$d = $word.ActiveDocument

function Add-StatusMarker($Marker, $Color, $SearchText, $AddTrailingSpace) {

    $target = $d.Content
    $found = $target.Find.Execute($SearchText, $true, $false, $false, $false,
                                   $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    $start = $target.Start
    $end = $target.End

    # Insert the marker text before the item's own text.
    $ins = $d.Range($start, $start)
    $ins.InsertBefore($Marker)

    $markerRange = $d.Range($start, $start + $Marker.Length)
    $markerRange.Font.Color = $Color
    $markerRange.Font.Name = "Arial"

    if ($AddTrailingSpace) {
        $textStart = $start + $Marker.Length
        $textEnd = $end + $Marker.Length
        $textRange = $d.Range($textStart, $textEnd)
        $textRange.Text = $textRange.Text + " "
    }
}

# "Worms detection algorithm:" sub-items
Add-StatusMarker "DONE" 52224 "Reading both image channels (5%)" $true

Add-StatusMarker "DONE" 52224 "Image segmentation (15%)" $true

Add-StatusMarker "PARTLY DONE" 39423 "Worms detection: background separation, objects labeling and counting (20%)" $false

# "Evaluation" sub-item
Add-StatusMarker "DONE" 52224 "Comparison of your warms detection results with the" $false
